$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 ("I0") and J1 ("IF")
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the formatting (style) of the existing header cell H1 onto I1:J1
# so they match the other header cells (bold, centered, bordered).
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data rows 2-23 for the new columns I (I0) and J (IF)
$data = @(
    @(6, 6),
    @(8, 8),
    @(9, 9),
    @(4, 4),
    @(5, 6),
    @(6, 6),
    @(9, 9),
    @(7, 8),
    @(8, 8),
    @(6, 7),
    @(7, 8),
    @(4, 5),
    @(8, 8),
    @(6, 6),
    @(8, 8),
    @(8, 8),
    @(8, 8),
    @(7, 7),
    @(4, 4),
    @(3, 3),
    @(8, 9),
    @(8, 8)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}
